# Update odds values on Sheet1 to match the 2024-11-13 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("G3").Value = 2.4
$ws.Range("Y3").Value = 9.5
$ws.Range("AB3").Value = 29
$ws.Range("AK3").Value = 34
$ws.Range("AP3").Value = 23
$ws.Range("AW3").Value = 5
$ws.Range("AY3").Value = 29

# Row 4
$ws.Range("J4").Value = 8
$ws.Range("M4").Value = 1.05
$ws.Range("N4").Value = 11
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 3.5
$ws.Range("Q4").Value = 1.95
$ws.Range("R4").Value = 1.9
$ws.Range("W4").Value = 9.5
$ws.Range("X4").Value = 41
$ws.Range("AJ4").Value = 9
$ws.Range("AK4").Value = 8.5
$ws.Range("AN4").Value = 9
$ws.Range("AQ4").Value = 201
$ws.Range("AW4").Value = 2.75

# Row 5
$ws.Range("AT5").Value = 1.83

# Row 7
$ws.Range("G7").Value = 1.87
$ws.Range("H7").Value = 3.2
$ws.Range("I7").Value = 4.1
$ws.Range("J7").Value = 2.4
$ws.Range("K7").Value = 2.07
$ws.Range("L7").Value = 4.45
$ws.Range("N7").Value = 6.5
$ws.Range("O7").Value = 1.34
$ws.Range("P7").Value = 2.75
$ws.Range("Q7").Value = 2
$ws.Range("R7").Value = 1.65
$ws.Range("S7").Value = 1.39
$ws.Range("T7").Value = 2.55
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.78
$ws.Range("W7").Value = 6.4
$ws.Range("Y7").Value = 8.25
$ws.Range("Z7").Value = 16
$ws.Range("AA7").Value = 15.5
$ws.Range("AB7").Value = 29
$ws.Range("AC7").Value = 8.25
$ws.Range("AD7").Value = 6.3
$ws.Range("AE7").Value = 15.5
$ws.Range("AF7").Value = 80
$ws.Range("AG7").Value = 700
$ws.Range("AH7").Value = 10.25
$ws.Range("AI7").Value = 22
$ws.Range("AJ7").Value = 14
$ws.Range("AK7").Value = 70
$ws.Range("AL7").Value = 45
$ws.Range("AN7").Value = 3.65
$ws.Range("AO7").Value = 9.25
$ws.Range("AP7").Value = 17.5
$ws.Range("AQ7").Value = 32
$ws.Range("AR7").Value = 60
$ws.Range("AS7").Value = 250
$ws.Range("AT7").Value = 2.52
$ws.Range("AU7").Value = 7.1
$ws.Range("AV7").Value = 65
$ws.Range("AW7").Value = 5.8
$ws.Range("AX7").Value = 23
$ws.Range("AY7").Value = 29
$ws.Range("AZ7").Value = 150
$ws.Range("BB7").Value = 400

# Row 8
$ws.Range("G8").Value = 2.02
$ws.Range("H8").Value = 3.35
$ws.Range("I8").Value = 3.4
$ws.Range("J8").Value = 2.6
$ws.Range("K8").Value = 2.1
$ws.Range("P8").Value = 3.25
$ws.Range("U8").Value = 1.62
$ws.Range("V8").Value = 2.02
$ws.Range("W8").Value = 8.25
$ws.Range("X8").Value = 10.25
$ws.Range("Z8").Value = 19
$ws.Range("AB8").Value = 24
$ws.Range("AC8").Value = 10.75
$ws.Range("AD8").Value = 6.5
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 55
$ws.Range("AH8").Value = 11.25
$ws.Range("AI8").Value = 19
$ws.Range("AK8").Value = 50
$ws.Range("AN8").Value = 3.95
$ws.Range("AO8").Value = 10.25
$ws.Range("AP8").Value = 18
$ws.Range("AQ8").Value = 37
$ws.Range("AT8").Value = 2.75
$ws.Range("AU8").Value = 6.8
$ws.Range("AY8").Value = 24
$ws.Range("AZ8").Value = 90
